# Add specifications for exterior blinds ("ext_blind") to the "properties" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("properties")

# New row right after the existing "ext_blind_test" row (row 12)
$ws.Range("A13").Value = "ext_blind"
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0

# Update selection to match the new active cell
$ws.Activate()
$ws.Range("E13").Select()
